$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record (Espinaca, Terminal La Palmera de La Serena) was
# added to the consolidated weekly log. It belongs right after the existing
# row for date 44643 (row 305) and before the row for date 44966 (old row
# 306), so insert a fresh row at position 306 — this pushes the old rows
# 306-450 down to 307-451.
$ws.Rows.Item(306).Insert()

# Excel's Insert() duplicates the formatting of the row above into the new
# blank row 306, but leaves the cell values empty. Seed it with the same
# "template" values (market/category/quality/price-band descriptors) as the
# record that used to sit at row 306 and now lives at row 307, then
# overwrite just the date (column D) and volume (column J) with the new
# record's own data.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(306, $col).Value = $ws.Cells.Item(307, $col).Value2
}

$ws.Cells.Item(306, 4).Value = 45134
$ws.Cells.Item(306, 10).Value = 1400
